$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "1.0-Huevos (unidad),3.0-Leche (litros),1.0-Vainilla (ml),2.0-Harina  (kg),"
$ws.Range("C3").Value = "5.0-Huevos (unidad),2.0-Harina  (kg),"
$ws.Range("C4").Value = "2.0-Huevos (unidad),1.0-Vainilla (ml),5.0-Harina  (kg),"
$ws.Range("C5").Value = "5.0-Huevos (unidad),5.0-Harina  (kg),"
$ws.Range("C6").Value = "5.0-Huevos (unidad),2.0-Limon (unidad),4.0-Harina  (kg),5.0-Crema (litros),"
$ws.Range("C7").Value = "2.0-Huevos (unidad),0.2-Leche (litros),0.1-Vainilla (ml),0.3-Harina  (kg),"
